$wb = $excel.ActiveWorkbook

# ----- Sheet: LP1912 -----
$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Value = 'Última actualización: 13:24:27'
$ws.Range("A3").Value = 'Total filas: 274'
$ws.Range("A124").Value = '08:33:47'
$ws.Range("C124").Value = '23_HERNANDEZ'
$ws.Range("D124").Value = 68
$ws.Range("A125").Value = '08:16:48'
$ws.Range("C125").Value = '215C_EL PATO'
$ws.Range("D125").Value = 85
$ws.Range("C229").Value = '17_ROMERO'
$ws.Range("C230").Value = '16_SANTA ANA'
$ws.Range("A245").Value = '13:24:27'
$ws.Range("B245").Value = '13:24'
$ws.Range("C245").Value = '17_ROMERO'
$ws.Range("D245").Value = 0
$ws.Range("A246").Value = '13:24:27'
$ws.Range("B246").Value = '13:25'
$ws.Range("C246").Value = '16_SANTA ANA'
$ws.Range("D246").Value = 1
$ws.Range("A247").Value = '11:43:19'
$ws.Range("B247").Value = '13:26'
$ws.Range("D247").Value = 103
$ws.Range("A248").Value = '11:43:19'
$ws.Range("B248").Value = '13:26'
$ws.Range("C248").Value = '15_ABASTO'
$ws.Range("D248").Value = 103
$ws.Range("A249").Value = '12:24:14'
$ws.Range("B249").Value = '13:27'
$ws.Range("C249").Value = '14_ABASTO'
$ws.Range("D249").Value = 63
$ws.Range("A250").Value = '11:58:46'
$ws.Range("B250").Value = '13:29'
$ws.Range("C250").Value = '17_ROMERO'
$ws.Range("D250").Value = 91
$ws.Range("A251").Value = '12:44:05'
$ws.Range("B251").Value = '13:32'
$ws.Range("C251").Value = '10_OLMOS'
$ws.Range("D251").Value = 48
$ws.Range("A252").Value = '12:57:33'
$ws.Range("B252").Value = '13:34'
$ws.Range("D252").Value = 37
$ws.Range("A253").Value = '13:24:27'
$ws.Range("B253").Value = '13:36'
$ws.Range("C253").Value = '15_ABASTO'
$ws.Range("D253").Value = 12
$ws.Range("A254").Value = '11:43:19'
$ws.Range("B254").Value = '13:37'
$ws.Range("C254").Value = '11_ETCHEVERRY'
$ws.Range("D254").Value = 114
$ws.Range("B255").Value = '13:40'
$ws.Range("C255").Value = '23_HERNANDEZ'
$ws.Range("D255").Value = 76
$ws.Range("A256").Value = '12:44:05'
$ws.Range("B256").Value = '13:41'
$ws.Range("C256").Value = '23_HERNANDEZ'
$ws.Range("D256").Value = 57
$ws.Range("A257").Value = '11:58:46'
$ws.Range("B257").Value = '13:46'
$ws.Range("C257").Value = '17_ROMERO'
$ws.Range("D257").Value = 108
$ws.Range("B258").Value = '13:47'
$ws.Range("C258").Value = '17_ROMERO'
$ws.Range("D258").Value = 83
$ws.Range("A259").Value = '12:57:33'
$ws.Range("B259").Value = '13:50'
$ws.Range("C259").Value = '11_ETCHEVERRY'
$ws.Range("D259").Value = 53
$ws.Range("B260").Value = '13:50'
$ws.Range("C260").Value = '215A_EL PATO'
$ws.Range("D260").Value = 112
$ws.Range("B261").Value = '13:51'
$ws.Range("C261").Value = '215A_EL PATO'
$ws.Range("D261").Value = 87
$ws.Range("A262").Value = '11:58:46'
$ws.Range("B262").Value = '13:56'
$ws.Range("C262").Value = '225_GOMEZ'
$ws.Range("D262").Value = 118
$ws.Range("A263").Value = '11:58:46'
$ws.Range("B263").Value = '13:56'
$ws.Range("C263").Value = '16_P MOR-167 Y 521'
$ws.Range("D263").Value = 118
$ws.Range("A264").Value = '12:24:14'
$ws.Range("B264").Value = '13:57'
$ws.Range("C264").Value = '16_P MOR-167 Y 521'
$ws.Range("D264").Value = 93
$ws.Range("A265").Value = '12:44:05'
$ws.Range("B265").Value = '14:04'
$ws.Range("C265").Value = '17_ROMERO'
$ws.Range("D265").Value = 80
$ws.Range("A266").Value = '13:24:27'
$ws.Range("B266").Value = '14:04'
$ws.Range("C266").Value = '23_HERNANDEZ'
$ws.Range("D266").Value = 40
$ws.Range("A267").Value = '13:24:27'
$ws.Range("B267").Value = '14:05'
$ws.Range("C267").Value = '11_ETCHEVERRY'
$ws.Range("D267").Value = 41
$ws.Range("A268").Value = '12:44:05'
$ws.Range("B268").Value = '14:05'
$ws.Range("C268").Value = '23_HERNANDEZ'
$ws.Range("D268").Value = 81
$ws.Range("B269").Value = '14:16'
$ws.Range("C269").Value = '27_EL RETIRO'
$ws.Range("D269").Value = 79
$ws.Range("A270").Value = '12:24:14'
$ws.Range("B270").Value = '14:17'
$ws.Range("C270").Value = '27_EL RETIRO'
$ws.Range("D270").Value = 113
$ws.Range("E270").Value = 'LP1912'
$ws.Range("A271").Value = '12:24:14'
$ws.Range("B271").Value = '14:20'
$ws.Range("C271").Value = '215C_EL PATO'
$ws.Range("D271").Value = 116
$ws.Range("E271").Value = 'LP1912'
$ws.Range("A272").Value = '12:24:14'
$ws.Range("B272").Value = '14:21'
$ws.Range("C272").Value = '26_HERNANDEZ'
$ws.Range("D272").Value = 117
$ws.Range("E272").Value = 'LP1912'
$ws.Range("A273").Value = '12:57:33'
$ws.Range("B273").Value = '14:45'
$ws.Range("C273").Value = '14_ABASTO'
$ws.Range("D273").Value = 108
$ws.Range("E273").Value = 'LP1912'
$ws.Range("A274").Value = '12:57:33'
$ws.Range("B274").Value = '14:56'
$ws.Range("C274").Value = '16_P MOR-SANTA ANA'
$ws.Range("D274").Value = 119
$ws.Range("E274").Value = 'LP1912'
$ws.Range("A275").Value = '13:24:27'
$ws.Range("B275").Value = '14:58'
$ws.Range("C275").Value = '215B_EL PATO'
$ws.Range("D275").Value = 94
$ws.Range("E275").Value = 'LP1912'
$ws.Range("A276").Value = '13:24:27'
$ws.Range("B276").Value = '15:00'
$ws.Range("C276").Value = '81_EL PELIGRO'
$ws.Range("D276").Value = 96
$ws.Range("E276").Value = 'LP1912'
$ws.Range("A277").Value = '13:24:27'
$ws.Range("B277").Value = '15:05'
$ws.Range("C277").Value = '10_OLMOS'
$ws.Range("D277").Value = 101
$ws.Range("E277").Value = 'LP1912'
$ws.Range("A278").Value = '13:24:27'
$ws.Range("B278").Value = '15:20'
$ws.Range("C278").Value = '15_ABASTO'
$ws.Range("D278").Value = 116
$ws.Range("E278").Value = 'LP1912'
$ws.Range("A279").Value = '13:24:27'
$ws.Range("B279").Value = '15:22'
$ws.Range("C279").Value = '26_HERNANDEZ'
$ws.Range("D279").Value = 118
$ws.Range("E279").Value = 'LP1912'

# ----- Sheet: LP1912-215 -----
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = 'Última actualización: 13:24:27'
$ws.Range("A3").Value = 'Total filas: 30'
$ws.Range("A35").Value = '13:24:27'
$ws.Range("B35").Value = '14:58'
$ws.Range("C35").Value = '215B_EL PATO'
$ws.Range("D35").Value = 94
$ws.Range("E35").Value = 'LP1912'

# ----- Sheet: 6203-6173 -----
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = 'Última actualización: 13:24:27'
